$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.127.80"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "2.470.54"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.53%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.513"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.335"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").Value = "2.908.85"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "67.085.45"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "2.441.62"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.57%  "
$ws.Range("E25").Value = "  -2.00%  "
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "2.599.26"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").Value = "0.0₃0899"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "498.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  +1.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "142.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.68%  "
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.582"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.21%  "
